$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A162").Value = "IMX-USD"
$ws.Range("A163").Value = "TAO-USD"
$ws.Range("A164").Value = "GRT-USD"
